$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.101.12"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.04%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.318.80"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.34%  "

# Row 4 - TetherUSD
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.96"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.19%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.37"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.89%  "

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.508"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.40%  "

# Row 8 - USDC
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.07%  "

# Row 9 - Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.517"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.63%  "

# Row 10 - Avalanche
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.94"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.27%  "

# Row 11 - Dogecoin
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.59%  "

# Row 12 - TRON
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.15%  "

# Row 13 - Chainlink
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.74"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.91%  "

# Row 14 - Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.89"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.49%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.679.20"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.21%  "

# Row 16 - WrappedEther
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.320.61"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.75%  "

# Row 17 - Polygon
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.11%  "

# Row 18 - WrappedBTC
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.009.40"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.18%  "

# Row 19 - InternetComputer(DFINITY)
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.12"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.48%  "

# Row 20 - Uniswap
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.85%  "

# Row 21 - ShibaInu
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.08%  "

# Row 22 - Litecoin
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.17"

# Row 23 - BitcoinCash
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.95"

# Row 24 - ImmutableX
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.31%  "

# Row 25 - PancakeSwap
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.45"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.62%  "

# Row 26 - Dai
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.04%  "

# Row 27 - EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.19"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.47%  "

# Row 28 - Monero
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.99"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.28%  "

# Row 29 - Cosmos
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.18"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.07%  "

# Row 31 - InjectiveProtocol
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.72"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.20%  "

# Row 32 - RenderToken
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.94"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +6.21%  "

# Row 33 - Filecoin
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.72%  "

# Row 34 - FirstDigitalUSD
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.04%  "

# Row 35 - Celestia
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.44"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +7.74%  "

# Row 36 - WEMIXToken
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.57%  "

# Row 37 - Hedera
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.79%  "

# Row 38 - ARBITRUM
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.45%  "

# Row 39 - Kaspa
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.30%  "

# Row 40 - LidoDAOToken
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.76"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.07%  "

# Row 41 - Stellar
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.27%  "

# Row 42 - Maker
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.995.51"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.34%  "

# Row 43 - VeChain
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0288"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.14%  "

# Row 44 - FraxShare (was ApeXProtocol)
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.18"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.47%  "

# Row 45 - ApeXProtocol (was FraxShare)
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.15"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -8.60%  "

# Row 46 - EnergySwap
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.42"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.78%  "

# Row 47 - NEARProtocol
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.84"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.84%  "

# Row 48 - BitcoinSV
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "76.43"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +8.68%  "

# Row 49 - MultiversX
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.89"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.78%  "

# Row 50 - RocketPoolETH
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.546.26"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.30%  "

# Row 51 - Stacks
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.03%  "
